$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Greater thann USD200 CR" / CHIPS rule row (row 11), shifting all
# subsequent rows up by one.
$ws.Rows(11).Delete()

# Move the active selection to the entire row 11 (mirrors the cursor position
# recorded after the edit in the source workbook).
$ws.Rows(11).Select()
